$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data occupies A2:D12 (row 1 holds header labels). The calibration data
# rows need to be re-sorted by time (column A) ascending, carrying the
# corresponding B/C/D curvature values along with each row.
$values = $ws.Range("A2:D12").Value2

$rows = New-Object System.Collections.ArrayList
for ($i = 1; $i -le 11; $i++) {
    $row = @($values[$i,1], $values[$i,2], $values[$i,3], $values[$i,4])
    $rows.Add($row) | Out-Null
}

$sorted = $rows | Sort-Object { $_[0] }

for ($i = 0; $i -lt $sorted.Count; $i++) {
    $r = 2 + $i
    $row = $sorted[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
}
